$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1848.9
$ws.Range("I40").Value = 1387.6666
$ws.Range("K40").Value = 1387.6666
$ws.Range("M40").Value = -1212.6666

$ws.Range("H64").Value = 4278
$ws.Range("I64").Value = 4000
$ws.Range("K64").Value = 4000
$ws.Range("M64").Value = -3752

$ws.Range("H67").Value = 4278
$ws.Range("I67").Value = 4000
$ws.Range("K67").Value = 4000
$ws.Range("M67").Value = -3142

$ws.Range("H132").Value = 1165.25
$ws.Range("I132").Value = 1131.7142
$ws.Range("K132").Value = 3395.1426
$ws.Range("M132").Value = -865.1425999999997

$ws.Range("H137").Value = 3540.1
$ws.Range("I137").Value = 3540.1
$ws.Range("K137").Value = 10620.3
$ws.Range("M137").Value = -8070.299999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3764.7556
$ws.Range("I32").Value = 1449.317
$ws.Range("K32").Value = 1449.317
$ws.Range("M32").Value = -1162.317

$ws.Range("H61").Value = 2698.2144
$ws.Range("I61").Value = 2554.5
$ws.Range("J61").Value = 3057.5
$ws.Range("K61").Value = 2554.5
$ws.Range("L61").Value = 3057.5
$ws.Range("M61").Value = -2342.5
$ws.Range("N61").Value = -3481.5

$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20496

$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21716

$ws.Range("H74").Value = 2657.8667
$ws.Range("I74").Value = 2490.5715
$ws.Range("K74").Value = 2490.5715
$ws.Range("M74").Value = -1616.5715

$ws.Range("H77").Value = 2657.8667
$ws.Range("I77").Value = 2490.5715
$ws.Range("K77").Value = 12452.8575
$ws.Range("M77").Value = -8084.8575

$ws.Range("H88").Value = 2568.7856
$ws.Range("I88").Value = 602.75
$ws.Range("K88").Value = 602.75
$ws.Range("M88").Value = -196.75

$ws.Range("H91").Value = 2568.7856
$ws.Range("I91").Value = 602.75
$ws.Range("K91").Value = 602.75
$ws.Range("M91").Value = 801.25

$ws.Range("H122").Value = 2723.5833
$ws.Range("I122").Value = 2723.5833
$ws.Range("K122").Value = 8170.749899999999
$ws.Range("M122").Value = -5720.749899999999

$ws.Range("H132").Value = 884.7692
$ws.Range("I132").Value = 884.7692
$ws.Range("K132").Value = 2654.3076
$ws.Range("M132").Value = -124.3076000000001

$ws.Range("H136").Value = 2698.2144
$ws.Range("I136").Value = 2554.5
$ws.Range("J136").Value = 3057.5
$ws.Range("K136").Value = 7663.5
$ws.Range("L136").Value = 9172.5
$ws.Range("M136").Value = -5113.5
$ws.Range("N136").Value = -14272.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1702.7858
$ws.Range("I105").Value = 1593.2222
$ws.Range("K105").Value = 1593.2222
$ws.Range("M105").Value = 153.7778000000001

$ws.Range("H134").Value = 5625.375
$ws.Range("I134").Value = 5914.2104
$ws.Range("K134").Value = 17742.6312
$ws.Range("M134").Value = -15207.6312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1672.8334
$ws.Range("I31").Value = 1367.5
$ws.Range("K31").Value = 1367.5
$ws.Range("M31").Value = -1072.5

$ws.Range("H34").Value = 1672.8334
$ws.Range("I34").Value = 1367.5
$ws.Range("K34").Value = 1367.5
$ws.Range("M34").Value = -1165.5

$ws.Range("H58").Value = 2217.111
$ws.Range("I58").Value = 2138.7144
$ws.Range("J58").Value = 2491.5
$ws.Range("K58").Value = 2138.7144
$ws.Range("L58").Value = 2491.5
$ws.Range("M58").Value = -1935.7144
$ws.Range("N58").Value = -2897.5

$ws.Range("H105").Value = 606
$ws.Range("I105").Value = 598.25
$ws.Range("K105").Value = 598.25
$ws.Range("M105").Value = 1148.75

$ws.Range("H132").Value = 5173.5557
$ws.Range("J132").Value = 5077.1665
$ws.Range("L132").Value = 15231.4995
$ws.Range("N132").Value = -20291.4995

$ws.Range("H136").Value = 2217.111
$ws.Range("I136").Value = 2138.7144
$ws.Range("J136").Value = 2491.5
$ws.Range("K136").Value = 6416.1432
$ws.Range("L136").Value = 7474.5
$ws.Range("M136").Value = -3866.1432
$ws.Range("N136").Value = -12574.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1058.4
$ws.Range("I5").Value = 1164.3334
$ws.Range("J5").Value = 899.5
$ws.Range("K5").Value = 3493.0002
$ws.Range("L5").Value = 2698.5
$ws.Range("M5").Value = -3381.0002
$ws.Range("N5").Value = -2922.5

$ws.Range("H34").Value = 1431.7142
$ws.Range("I34").Value = 600.2
$ws.Range("K34").Value = 1800.6
$ws.Range("M34").Value = -1716.6

$ws.Range("H135").Value = 1058.4
$ws.Range("I135").Value = 1164.3334
$ws.Range("J135").Value = 899.5
$ws.Range("K135").Value = 10479.0006
$ws.Range("L135").Value = 8095.5
$ws.Range("M135").Value = -7944.000599999999
$ws.Range("N135").Value = -13165.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 55
$ws.Range("I5").Value = 55
$ws.Range("K5").Value = 55
$ws.Range("M5").Value = 57

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws.Range("H126").Value = 7330.1
$ws.Range("I126").Value = 3662.625
$ws.Range("J126").Value = 22000
$ws.Range("K126").Value = 10987.875
$ws.Range("L126").Value = 66000
$ws.Range("M126").Value = -8517.875
$ws.Range("N126").Value = -70940

$ws.Range("H132").Value = 2248.25
$ws.Range("I132").Value = 2331
$ws.Range("K132").Value = 6993
$ws.Range("M132").Value = -4463

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 17252.818
$ws.Range("I16").Value = 3848.625
$ws.Range("J16").Value = 52997.332
$ws.Range("K16").Value = 3848.625
$ws.Range("L16").Value = 52997.332
$ws.Range("M16").Value = -3678.625
$ws.Range("N16").Value = -53337.332

$ws.Range("H46").Value = 2266.3333
$ws.Range("I46").Value = 2400
$ws.Range("J46").Value = 1999
$ws.Range("K46").Value = 2400
$ws.Range("L46").Value = 1999
$ws.Range("M46").Value = -2212
$ws.Range("N46").Value = -2375

$ws.Range("H93").Value = 2350.6667
$ws.Range("I93").Value = 2951.5
$ws.Range("J93").Value = 1149
$ws.Range("K93").Value = 2951.5
$ws.Range("L93").Value = 1149
$ws.Range("M93").Value = -1703.5
$ws.Range("N93").Value = -3645

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 48423.668
$ws.Range("J68").Value = 48423.668
$ws.Range("L68").Value = 48423.668
$ws.Range("N68").Value = -50045.668

$ws.Range("H71").Value = 48423.668
$ws.Range("J71").Value = 48423.668
$ws.Range("L71").Value = 145271.004
$ws.Range("N71").Value = -153383.004

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H132").Value = 3713.8333
$ws.Range("I132").Value = 1001
$ws.Range("K132").Value = 3003
$ws.Range("M132").Value = -473

$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120

$ws.Range("H136").Value = 8297.538
$ws.Range("I136").Value = 7216.3335
$ws.Range("K136").Value = 21649.0005
$ws.Range("M136").Value = -19099.0005
